# Updates scraped "profit" figures (currentAveragePrice / NQ / HQ / LevePrice / LeveProfit
# columns H-N) across several leve sheets, per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(129, 8).Value = 870.2373
$ws.Cells.Item(129, 9).Value = 386.16666
$ws.Cells.Item(129, 10).Value = 993.8298
$ws.Cells.Item(129, 11).Value = 1158.49998
$ws.Cells.Item(129, 12).Value = 2981.4894
$ws.Cells.Item(129, 13).Value = 3841.50002
$ws.Cells.Item(129, 14).Value = -12981.4894

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 4067.6924
$ws.Cells.Item(63, 9).Value = 2976
$ws.Cells.Item(63, 11).Value = 2976
$ws.Cells.Item(63, 13).Value = -2290

$ws.Cells.Item(66, 8).Value = 4067.6924
$ws.Cells.Item(66, 9).Value = 2976
$ws.Cells.Item(66, 11).Value = 14880
$ws.Cells.Item(66, 13).Value = -11448

$ws.Cells.Item(88, 8).Value = 18658.834
$ws.Cells.Item(88, 9).Value = 1999.6666
$ws.Cells.Item(88, 10).Value = 35318
$ws.Cells.Item(88, 11).Value = 1999.6666
$ws.Cells.Item(88, 12).Value = 35318
$ws.Cells.Item(88, 13).Value = -1593.6666
$ws.Cells.Item(88, 14).Value = -36130

$ws.Cells.Item(91, 8).Value = 18658.834
$ws.Cells.Item(91, 9).Value = 1999.6666
$ws.Cells.Item(91, 10).Value = 35318
$ws.Cells.Item(91, 11).Value = 1999.6666
$ws.Cells.Item(91, 12).Value = 35318
$ws.Cells.Item(91, 13).Value = -595.6666
$ws.Cells.Item(91, 14).Value = -38126

$ws.Cells.Item(110, 8).Value = 1349.1786
$ws.Cells.Item(110, 9).Value = 1575.5714
$ws.Cells.Item(110, 10).Value = 670
$ws.Cells.Item(110, 11).Value = 1575.5714
$ws.Cells.Item(110, 12).Value = 670
$ws.Cells.Item(110, 13).Value = 469.4286
$ws.Cells.Item(110, 14).Value = -4760

$ws.Cells.Item(132, 8).Value = 2123.8
$ws.Cells.Item(132, 9).Value = 1770.0278
$ws.Cells.Item(132, 10).Value = 3033.5
$ws.Cells.Item(132, 11).Value = 5310.0834
$ws.Cells.Item(132, 12).Value = 9100.5
$ws.Cells.Item(132, 13).Value = -2780.0834
$ws.Cells.Item(132, 14).Value = -14160.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(35, 8).Value = 31291.334
$ws.Cells.Item(35, 10).Value = 31291.334
$ws.Cells.Item(35, 12).Value = 31291.334
$ws.Cells.Item(35, 14).Value = -31911.334

$ws.Cells.Item(82, 8).Value = 34163
$ws.Cells.Item(82, 10).Value = 35197.6
$ws.Cells.Item(82, 12).Value = 35197.6
$ws.Cells.Item(82, 14).Value = -35963.6

$ws.Cells.Item(85, 8).Value = 34163
$ws.Cells.Item(85, 10).Value = 35197.6
$ws.Cells.Item(85, 12).Value = 35197.6
$ws.Cells.Item(85, 14).Value = -37849.6

$ws.Cells.Item(86, 8).Value = 2017.9131
$ws.Cells.Item(86, 9).Value = 2182.5833
$ws.Cells.Item(86, 10).Value = 1838.2727
$ws.Cells.Item(86, 11).Value = 2182.5833
$ws.Cells.Item(86, 12).Value = 1838.2727
$ws.Cells.Item(86, 13).Value = -1059.5833
$ws.Cells.Item(86, 14).Value = -4084.2727

$ws.Cells.Item(89, 8).Value = 2017.9131
$ws.Cells.Item(89, 9).Value = 2182.5833
$ws.Cells.Item(89, 10).Value = 1838.2727
$ws.Cells.Item(89, 11).Value = 10912.9165
$ws.Cells.Item(89, 12).Value = 9191.363499999999
$ws.Cells.Item(89, 13).Value = -5296.916499999999
$ws.Cells.Item(89, 14).Value = -20423.3635

$ws.Cells.Item(107, 8).Value = 6230.048
$ws.Cells.Item(107, 9).Value = 6248.9443
$ws.Cells.Item(107, 10).Value = 6116.6665
$ws.Cells.Item(107, 11).Value = 6248.9443
$ws.Cells.Item(107, 12).Value = 6116.6665
$ws.Cells.Item(107, 13).Value = -4328.9443
$ws.Cells.Item(107, 14).Value = -9956.666499999999

$ws.Cells.Item(113, 8).Value = 4940
$ws.Cells.Item(113, 9).Value = 4940
$ws.Cells.Item(113, 11).Value = 4940
$ws.Cells.Item(113, 13).Value = -2770

$ws.Cells.Item(122, 8).Value = 45511.25
$ws.Cells.Item(122, 10).Value = 45511.25
$ws.Cells.Item(122, 12).Value = 45511.25
$ws.Cells.Item(122, 14).Value = -55311.25

$ws.Cells.Item(134, 8).Value = 2350.739
$ws.Cells.Item(134, 9).Value = 1761.6
$ws.Cells.Item(134, 10).Value = 3455.375
$ws.Cells.Item(134, 11).Value = 5284.799999999999
$ws.Cells.Item(134, 12).Value = 10366.125
$ws.Cells.Item(134, 13).Value = -2749.799999999999
$ws.Cells.Item(134, 14).Value = -15436.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(50, 8).Value = 9129.333000000001
$ws.Cells.Item(50, 10).Value = 9129.333000000001
$ws.Cells.Item(50, 12).Value = 9129.333000000001
$ws.Cells.Item(50, 14).Value = -10379.333

$ws.Cells.Item(51, 8).Value = 9067
$ws.Cells.Item(51, 10).Value = 9067
$ws.Cells.Item(51, 12).Value = 9067
$ws.Cells.Item(51, 14).Value = -10539

$ws.Cells.Item(61, 8).Value = 9067
$ws.Cells.Item(61, 10).Value = 9067
$ws.Cells.Item(61, 12).Value = 9067
$ws.Cells.Item(61, 14).Value = -9763

$ws.Cells.Item(68, 8).Value = 16979
$ws.Cells.Item(68, 10).Value = 16979
$ws.Cells.Item(68, 12).Value = 16979
$ws.Cells.Item(68, 14).Value = -18477

$ws.Cells.Item(71, 8).Value = 16979
$ws.Cells.Item(71, 10).Value = 16979
$ws.Cells.Item(71, 12).Value = 50937
$ws.Cells.Item(71, 14).Value = -58425

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 570.8889
$ws.Cells.Item(5, 9).Value = 378.53125
$ws.Cells.Item(5, 10).Value = 1044.3846
$ws.Cells.Item(5, 11).Value = 1135.59375
$ws.Cells.Item(5, 12).Value = 3133.1538
$ws.Cells.Item(5, 13).Value = -1023.59375
$ws.Cells.Item(5, 14).Value = -3357.1538

$ws.Cells.Item(40, 8).Value = 253.33333
$ws.Cells.Item(40, 9).Value = 68.57143000000001
$ws.Cells.Item(40, 11).Value = 274.28572
$ws.Cells.Item(40, 13).Value = -205.28572

$ws.Cells.Item(131, 8).Value = 879.48
$ws.Cells.Item(131, 9).Value = 489.69232
$ws.Cells.Item(131, 10).Value = 937.7241
$ws.Cells.Item(131, 11).Value = 1469.07696
$ws.Cells.Item(131, 12).Value = 2813.1723
$ws.Cells.Item(131, 13).Value = 3570.92304
$ws.Cells.Item(131, 14).Value = -12893.1723

$ws.Cells.Item(132, 8).Value = 528658.3
$ws.Cells.Item(132, 9).Value = 1317049.8
$ws.Cells.Item(132, 10).Value = 3064
$ws.Cells.Item(132, 11).Value = 11853448.2
$ws.Cells.Item(132, 12).Value = 27576
$ws.Cells.Item(132, 13).Value = -11850918.2
$ws.Cells.Item(132, 14).Value = -32636

$ws.Cells.Item(135, 8).Value = 570.8889
$ws.Cells.Item(135, 9).Value = 378.53125
$ws.Cells.Item(135, 10).Value = 1044.3846
$ws.Cells.Item(135, 11).Value = 3406.78125
$ws.Cells.Item(135, 12).Value = 9399.4614
$ws.Cells.Item(135, 13).Value = -871.78125
$ws.Cells.Item(135, 14).Value = -14469.4614

$ws.Cells.Item(141, 8).Value = 4592.522
$ws.Cells.Item(141, 9).Value = 2601.75
$ws.Cells.Item(141, 10).Value = 9142.857
$ws.Cells.Item(141, 11).Value = 7805.25
$ws.Cells.Item(141, 12).Value = 27428.571
$ws.Cells.Item(141, 13).Value = -2625.25
$ws.Cells.Item(141, 14).Value = -37788.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 63.77778
$ws.Cells.Item(2, 9).Value = 48.4
$ws.Cells.Item(2, 10).Value = 83
$ws.Cells.Item(2, 11).Value = 48.4
$ws.Cells.Item(2, 12).Value = 83
$ws.Cells.Item(2, 13).Value = 64.59999999999999
$ws.Cells.Item(2, 14).Value = -309

$ws.Cells.Item(53, 8).Value = 6000
$ws.Cells.Item(53, 9).Value = 0
$ws.Cells.Item(53, 11).Value = 0
$ws.Cells.Item(53, 13).ClearContents()

$ws.Cells.Item(70, 8).Value = 35923.91
$ws.Cells.Item(70, 9).Value = 44315.383
$ws.Cells.Item(70, 10).Value = 4755.5713
$ws.Cells.Item(70, 11).Value = 44315.383
$ws.Cells.Item(70, 12).Value = 4755.5713
$ws.Cells.Item(70, 13).Value = -44045.383
$ws.Cells.Item(70, 14).Value = -5295.5713

$ws.Cells.Item(73, 8).Value = 35923.91
$ws.Cells.Item(73, 9).Value = 44315.383
$ws.Cells.Item(73, 10).Value = 4755.5713
$ws.Cells.Item(73, 11).Value = 44315.383
$ws.Cells.Item(73, 12).Value = 4755.5713
$ws.Cells.Item(73, 13).Value = -43379.383
$ws.Cells.Item(73, 14).Value = -6627.5713

$ws.Cells.Item(97, 8).Value = 1784.5454
$ws.Cells.Item(97, 9).Value = 1763
$ws.Cells.Item(97, 11).Value = 1763
$ws.Cells.Item(97, 13).Value = -1267

$ws.Cells.Item(102, 8).Value = 1007.3333
$ws.Cells.Item(102, 9).Value = 855.6667
$ws.Cells.Item(102, 11).Value = 855.6667
$ws.Cells.Item(102, 13).Value = 766.3333

$ws.Cells.Item(132, 8).Value = 2235.366
$ws.Cells.Item(132, 9).Value = 2059
$ws.Cells.Item(132, 10).Value = 2661.5833
$ws.Cells.Item(132, 11).Value = 6177
$ws.Cells.Item(132, 12).Value = 7984.749899999999
$ws.Cells.Item(132, 13).Value = -3647
$ws.Cells.Item(132, 14).Value = -13044.7499

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 656
$ws.Cells.Item(55, 9).Value = 114.61539
$ws.Cells.Item(55, 10).Value = 1095.875
$ws.Cells.Item(55, 11).Value = 114.61539
$ws.Cells.Item(55, 12).Value = 1095.875
$ws.Cells.Item(55, 13).Value = 58.38461
$ws.Cells.Item(55, 14).Value = -1441.875

$ws.Cells.Item(93, 8).Value = 1095.093
$ws.Cells.Item(93, 9).Value = 1021.36365
$ws.Cells.Item(93, 10).Value = 1338.4
$ws.Cells.Item(93, 11).Value = 1021.36365
$ws.Cells.Item(93, 12).Value = 1338.4
$ws.Cells.Item(93, 13).Value = 226.63635
$ws.Cells.Item(93, 14).Value = -3834.4

$ws.Cells.Item(133, 8).Value = 65580
$ws.Cells.Item(133, 10).Value = 65580
$ws.Cells.Item(133, 12).Value = 65580
$ws.Cells.Item(133, 14).Value = -70640

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(109, 8).Value = 19441.334
$ws.Cells.Item(109, 10).Value = 19441.334
$ws.Cells.Item(109, 12).Value = 19441.334
$ws.Cells.Item(109, 14).Value = -22215.334

$ws.Cells.Item(132, 8).Value = 1873.8776
$ws.Cells.Item(132, 9).Value = 1946.8611
$ws.Cells.Item(132, 10).Value = 1671.7693
$ws.Cells.Item(132, 11).Value = 5840.5833
$ws.Cells.Item(132, 12).Value = 5015.3079
$ws.Cells.Item(132, 13).Value = -3310.5833
$ws.Cells.Item(132, 14).Value = -10075.3079

Write-Host "Sheets updated: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR"
